{"js": "// The edit replaces the header date and every three-digit x one-digit\n// multiplication \"problem=answer\" string in the practice-sheet table with\n// new values. Every old string is unique within the document, so each\n// pair can be resolved with a single Word body search-and-replace; this\n// preserves paragraph/run structure and formatting (fonts, size, etc.)\n// because only the <w:t> text content is touched.\nconst replacements = [\n  [\"2026-01-16 Friday\", \"2026-01-17 Saturday\"],\n  [\"660\u00d75=3300\", \"846\u00d76=5076\"],\n  [\"784\u00d73=2352\", \"878\u00d74=3512\"],\n  [\"667\u00d78=5336\", \"193\u00d72=386\"],\n  [\"269\u00d75=1345\", \"444\u00d74=1776\"],\n  [\"197\u00d74=788\", \"823\u00d77=5761\"],\n  [\"282\u00d73=846\", \"759\u00d77=5313\"],\n  [\"212\u00d72=424\", \"435\u00d74=1740\"],\n  [\"868\u00d76=5208\", \"461\u00d74=1844\"],\n  [\"201\u00d74=804\", \"994\u00d74=3976\"],\n  [\"254\u00d74=1016\", \"286\u00d75=1430\"],\n  [\"818\u00d78=6544\", \"173\u00d78=1384\"],\n  [\"387\u00d73=1161\", \"937\u00d72=1874\"],\n  [\"603\u00d74=2412\", \"378\u00d73=1134\"],\n  [\"499\u00d76=2994\", \"633\u00d78=5064\"],\n  [\"896\u00d72=1792\", \"186\u00d74=744\"],\n  [\"695\u00d73=2085\", \"824\u00d77=5768\"],\n  [\"168\u00d79=1512\", \"516\u00d73=1548\"],\n  [\"445\u00d75=2225\", \"531\u00d78=4248\"],\n  [\"879\u00d76=5274\", \"814\u00d73=2442\"],\n  [\"494\u00d76=2964\", \"967\u00d78=7736\"],\n  [\"736\u00d78=5888\", \"475\u00d72=950\"],\n  [\"187\u00d74=748\", \"420\u00d74=1680\"],\n  [\"813\u00d74=3252\", \"165\u00d77=1155\"],\n  [\"389\u00d72=778\", \"988\u00d72=1976\"],\n  [\"632\u00d74=2528\", \"678\u00d74=2712\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and all 25 three-digit x one-digit multiplication\n# answers in the practice-sheet table. Each old value is unique in the\n# document, so Find/Replace (wdReplaceAll) on each pair is unambiguous and\n# preserves the original run formatting (fonts/size) of the matched text.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-01-16 Friday\", \"2026-01-17 Saturday\"),\n    @(\"660\u00d75=3300\", \"846\u00d76=5076\"),\n    @(\"784\u00d73=2352\", \"878\u00d74=3512\"),\n    @(\"667\u00d78=5336\", \"193\u00d72=386\"),\n    @(\"269\u00d75=1345\", \"444\u00d74=1776\"),\n    @(\"197\u00d74=788\", \"823\u00d77=5761\"),\n    @(\"282\u00d73=846\", \"759\u00d77=5313\"),\n    @(\"212\u00d72=424\", \"435\u00d74=1740\"),\n    @(\"868\u00d76=5208\", \"461\u00d74=1844\"),\n    @(\"201\u00d74=804\", \"994\u00d74=3976\"),\n    @(\"254\u00d74=1016\", \"286\u00d75=1430\"),\n    @(\"818\u00d78=6544\", \"173\u00d78=1384\"),\n    @(\"387\u00d73=1161\", \"937\u00d72=1874\"),\n    @(\"603\u00d74=2412\", \"378\u00d73=1134\"),\n    @(\"499\u00d76=2994\", \"633\u00d78=5064\"),\n    @(\"896\u00d72=1792\", \"186\u00d74=744\"),\n    @(\"695\u00d73=2085\", \"824\u00d77=5768\"),\n    @(\"168\u00d79=1512\", \"516\u00d73=1548\"),\n    @(\"445\u00d75=2225\", \"531\u00d78=4248\"),\n    @(\"879\u00d76=5274\", \"814\u00d73=2442\"),\n    @(\"494\u00d76=2964\", \"967\u00d78=7736\"),\n    @(\"736\u00d78=5888\", \"475\u00d72=950\"),\n    @(\"187\u00d74=748\", \"420\u00d74=1680\"),\n    @(\"813\u00d74=3252\", \"165\u00d77=1155\"),\n    @(\"389\u00d72=778\", \"988\u00d72=1976\"),\n    @(\"632\u00d74=2528\", \"678\u00d74=2712\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}\n"}
